# Sprint2 planning and closing Sprint1
# Mark the three still-"In progress" Sprint1 tasks (Document implementation /
# Timeline / Document implementation) as Done, logging their remaining effort
# against Day 14 so the burndown totals close out at zero, and move the
# worksheet selection to F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

# Rows 8, 11 and 13 were still "In progress" with remaining effort booked on
# Day 14 (column T) to finish the task out.
$ws.Range("F8").Value  = "Done"
$ws.Range("T8").Value  = 2

$ws.Range("F11").Value = "Done"
$ws.Range("T11").Value = 4

$ws.Range("F13").Value = "Done"
$ws.Range("T13").Value = 2

# Leave the cursor where the author left it when saving.
$ws.Range("F11").Select() | Out-Null
